$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set NumberFormat to Text for all changed cells first to avoid Excel
# auto-converting numeric-looking strings (e.g. "46.90" -> 46.9, "111.00" -> 111)
# or values that could be parsed as dates.
$changedCells = @("D2","E2","D3","E3","E4","D5","E5","E6","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D27","E27","D28","E28","D29","E29","D30","E30","D31","E31","D32","E32","D33","E33","D34","E34","D35","E35","D36","E36","D37","E37","D38","E38","D39","E39","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","B46","C46","D46","E46","B47","C47","D47","E47","D48","E48","D49","E49","D50","E50","E51")
foreach ($addr in $changedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.798.57"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "1.850.14"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "335.51"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.3859"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "46.90"
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("D10").Value = "0.07924"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").Value = "0.9690"
$ws.Range("E11").Value = "  -2.89%  "
$ws.Range("D12").Value = "21.32"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "1.851.08"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "5.880"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "7.132"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "90.93"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("D18").Value = "0.06619"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "0.00001029"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").Value = "17.28"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "1.008"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "27.794.01"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("D23").Value = "5.346"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "10.80"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").Value = "2.289"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").Value = "2.072.93"
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "159.28"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "19.46"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "2.068"
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("D30").Value = "5.389"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").Value = "118.59"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "0.09441"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").Value = "0.9412"
$ws.Range("E33").Value = "  -3.55%  "
$ws.Range("D34").Value = "3.594"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "5.260"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").Value = "1.328"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "0.06025"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "0.02209"
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("D39").Value = "8.220"
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").Value = "1.162"
$ws.Range("E41").Value = "  -1.53%  "
$ws.Range("D42").Value = "0.5812"
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("D43").Value = "0.1848"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").Value = "10.08"
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("D45").Value = "1.284"
$ws.Range("E45").Value = "  +3.55%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5457"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "11.93"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").Value = "1.937"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("D49").Value = "0.06838"
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("D50").Value = "111.00"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("E51").Value = "  -32.37%  "
